# Commit: Tue, Jul 21, 2020 10:05:35 PM
#
# 1) Slide 6's finance-sources table switches to a different built-in
#    PowerPoint table style.
# 2) The presentation's applied design/theme colours change from the
#    "Integral" palette to the standard "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{5ABBCF13-7F18-42B8-B3F4-81D4AAC5D157}")

# --- 2. Swap the design's colour scheme to the Office Theme palette -----
$master = $p.SlideMaster

$officeThemeColors = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $master.ColorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
